# Generate Report for Handoff
#
# The localization-status report is regenerated: the two tracked files
# (1e7a8862-... and a73e1806-...) swap row order on every sheet, and the
# 1e7a8862 file has now moved from "In Translation" to "Ready for handoff"
# with a fresh handoff timestamp, while a73e1806 is still "In Translation".

$wb = $excel.ActiveWorkbook

function Set-HyperlinkDisplay {
    param($ws, [string]$addr, [string]$text)
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $h.TextToDisplay = $text
        }
    }
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "a73e1806-9ed2-4002-bf8d-6a5e8294ef4f.md"
$ov.Range("B2").Value = "In Translation"
$ov.Range("C2").Value = "In Translation"
$ov.Range("D2").Value = "2016-03-21 22:17:51"

$ov.Range("A3").Value = "1e7a8862-0e2f-4686-b980-af3634e00169.md"
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-03-21 22:18:41"

Set-HyperlinkDisplay $ov '$A$2' "a73e1806-9ed2-4002-bf8d-6a5e8294ef4f.md"
Set-HyperlinkDisplay $ov '$A$3' "1e7a8862-0e2f-4686-b980-af3634e00169.md"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "a73e1806-9ed2-4002-bf8d-6a5e8294ef4f.md"
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = "In Translation"
$zh.Range("D2").Value = "a73e1806-9ed2-4002-bf8d-6a5e8294ef4f.8032fef6cf08f9123ca86f621da026d3d1029544.zh-cn.xlf"
$zh.Range("E2").Value = "2016-03-21 22:17:48"
$zh.Range("H2").Value = "0001-01-01 00:00:00"
$zh.Range("J2").Value = "Include"

$zh.Range("A3").Value = "1e7a8862-0e2f-4686-b980-af3634e00169.md"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("D3").Value = "1e7a8862-0e2f-4686-b980-af3634e00169.2826213ef75e7ab619c75fa2c334e542a0d4adba.zh-cn.xlf"
$zh.Range("E3").Value = "2016-03-21 22:18:37"
$zh.Range("H3").Value = "0001-01-01 00:00:00"
$zh.Range("J3").Value = "Include"

Set-HyperlinkDisplay $zh '$A$2' "a73e1806-9ed2-4002-bf8d-6a5e8294ef4f.md"
Set-HyperlinkDisplay $zh '$D$2' "a73e1806-9ed2-4002-bf8d-6a5e8294ef4f.8032fef6cf08f9123ca86f621da026d3d1029544.zh-cn.xlf"
Set-HyperlinkDisplay $zh '$A$3' "1e7a8862-0e2f-4686-b980-af3634e00169.md"
Set-HyperlinkDisplay $zh '$D$3' "1e7a8862-0e2f-4686-b980-af3634e00169.2826213ef75e7ab619c75fa2c334e542a0d4adba.zh-cn.xlf"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "a73e1806-9ed2-4002-bf8d-6a5e8294ef4f.md"
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = "In Translation"
$de.Range("D2").Value = "a73e1806-9ed2-4002-bf8d-6a5e8294ef4f.8032fef6cf08f9123ca86f621da026d3d1029544.de-de.xlf"
$de.Range("E2").Value = "2016-03-21 22:17:51"
$de.Range("H2").Value = "0001-01-01 00:00:00"
$de.Range("J2").Value = "Include"

$de.Range("A3").Value = "1e7a8862-0e2f-4686-b980-af3634e00169.md"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("D3").Value = "1e7a8862-0e2f-4686-b980-af3634e00169.2826213ef75e7ab619c75fa2c334e542a0d4adba.de-de.xlf"
$de.Range("E3").Value = "2016-03-21 22:18:41"
$de.Range("H3").Value = "0001-01-01 00:00:00"
$de.Range("J3").Value = "Include"

Set-HyperlinkDisplay $de '$A$2' "a73e1806-9ed2-4002-bf8d-6a5e8294ef4f.md"
Set-HyperlinkDisplay $de '$D$2' "a73e1806-9ed2-4002-bf8d-6a5e8294ef4f.8032fef6cf08f9123ca86f621da026d3d1029544.de-de.xlf"
Set-HyperlinkDisplay $de '$A$3' "1e7a8862-0e2f-4686-b980-af3634e00169.md"
Set-HyperlinkDisplay $de '$D$3' "1e7a8862-0e2f-4686-b980-af3634e00169.2826213ef75e7ab619c75fa2c334e542a0d4adba.de-de.xlf"

Write-Output "Report regenerated for handoff"
